# Natmi following Dr Hou advice
#
# The ligand/receptor pair sheet (Il11 -> Il6st) is rebuilt to include a
# third "Sending cluster" / "Target cluster" category (ECs), in addition to
# the pre-existing FAPs and sCs clusters, giving a full 3x3 cross of
# sending/target clusters (9 data rows instead of 6) with refreshed
# specificity/expression statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# One row per (Sending cluster, Target cluster) combination, in the final
# on-sheet order: columns A:T = Sending cluster, Ligand symbol, Receptor
# symbol, Target cluster, then the 16 numeric expression/specificity columns.
$rows = @(
    @("ECs",  "Il11", "Il6st", "ECs",  1, 0.3333333333333333, 0.04156866666666666, 0.124706,  0.02383845226880681, 0.02383845226880681, 3, 1, 39.96608766666667,  119.898263, 0.2616165719423124, 0.2616165719423124, 1.661336976186445,  14.952032785678,    0.006236534162975677, 0.006236534162975677),
    @("ECs",  "Il11", "Il6st", "FAPs", 1, 0.3333333333333333, 0.04156866666666666, 0.124706,  0.02383845226880681, 0.02383845226880681, 3, 1, 85.11961100000001,  255.358833, 0.5571899111219771, 0.557189911121977,  3.538308736455333,  31.844778628098,    0.01328254510094196,  0.01328254510094196),
    @("ECs",  "Il11", "Il6st", "sCs",  1, 0.3333333333333333, 0.04156866666666666, 0.124706,  0.02383845226880681, 0.02383845226880681, 3, 1, 27.68018833333333,  83.040565,  0.1811935169357105, 0.1811935169357105, 1.150628522098889,  10.35565669889,     0.004319373004889172, 0.004319373004889172),
    @("FAPs", "Il11", "Il6st", "ECs",  3, 1,                  1.218041,            3.654123,  0.6985119939686074, 0.6985119939686074,  3, 1, 39.96608766666667,  119.898263, 0.2616165719423124, 0.2616165719423124, 48.68033338759435,  438.1230004883491,  0.1827423133226563,   0.1827423133226563),
    @("FAPs", "Il11", "Il6st", "FAPs", 3, 1,                  1.218041,            3.654123,  0.6985119939686074, 0.6985119939686074,  3, 1, 85.11961100000001,  255.358833, 0.5571899111219771, 0.557189911121977,  103.679176102051,   933.1125849184591,  0.3892038358370034,   0.3892038358370033),
    @("FAPs", "Il11", "Il6st", "sCs",  3, 1,                  1.218041,            3.654123,  0.6985119939686074, 0.6985119939686074,  3, 1, 27.68018833333333,  83.040565,  0.1811935169357105, 0.1811935169357105, 33.71560427772167,  303.440438499495,   0.1265658448089478,   0.1265658448089478),
    @("sCs",  "Il11", "Il6st", "ECs",  3, 1,                  0.4841556666666667,  1.452467,  0.2776495537625858, 0.2776495537625858,  3, 1, 39.96608766666667,  119.898263, 0.2616165719423124, 0.2616165719423124, 19.34980781831345,  174.148270364821,   0.07263772445668046,  0.07263772445668046),
    @("sCs",  "Il11", "Il6st", "FAPs", 3, 1,                  0.4841556666666667,  1.452467,  0.2776495537625858, 0.2776495537625858,  3, 1, 85.11961100000001,  255.358833, 0.5571899111219771, 0.557189911121977,  41.21114201011233,  370.900278091011,   0.1547035301840318,   0.1547035301840317),
    @("sCs",  "Il11", "Il6st", "sCs",  3, 1,                  0.4841556666666667,  1.452467,  0.2776495537625858, 0.2776495537625858,  3, 1, 27.68018833333333,  83.040565,  0.1811935169357105, 0.1811935169357105, 13.40152003598389,  120.613680323855,   0.05030829912187353,  0.05030829912187353)
)

$columns = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$rowIndex = 2
foreach ($rowData in $rows) {
    for ($c = 0; $c -lt $columns.Length; $c++) {
        $ws.Range($columns[$c] + $rowIndex).Value = $rowData[$c]
    }
    $rowIndex++
}
